$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Periodo Mora" (column E) labels for rows 16-22 - reorder from
# descending (2410..2404) to ascending (2404..2410).
$periods = @("2404","2405","2406","2407","2408","2409","2410")
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
}

# Update "Salario Basico" (column F) values - the partial-period salary
# value (53334) moves from the last row to the first row.
$salarios = @(53334, 80000, 80000, 80000, 80000, 80000, 80000)
for ($i = 0; $i -lt $salarios.Length; $i++) {
    $row = 16 + $i
    $ws.Range("F$row").Value = $salarios[$i]
}

# Update "Valor Mora" (column G) values for rows 16-22 from 0 to 2000000.
for ($row = 16; $row -le 22; $row++) {
    $ws.Range("G$row").Value = 2000000
}
